# Daily attendance processing - 2025-12-16 01:32:31
# Normalize the "Recorded By" (column G) entries: move the last
# comma-separated recorder to the front of the list, leaving the
# "admin@admin.com, System" combination untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val

    if ($text -eq "admin@admin.com, System") { continue }

    if ($text.Contains(",")) {
        $parts = $text -split ",\s*"
        if ($parts.Count -gt 1) {
            $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
            $newText = [string]::Join(", ", $rotated)
            $cell.Value = $newText
        }
    }
}
